$wb = $excel.ActiveWorkbook

# --- thule_lakepool_iso_2018vs2019_p (sheet10): rebuild attribute table with new/reordered rows ---
$ws = $wb.Worksheets.Item("thule_lakepool_iso_2018vs2019_p")

$ws.Range("A1").Value = 'attribute name'
$ws.Range("B1").Value = 'description'
$ws.Range("C1").Value = 'options'
$ws.Range("A2").Value = 'id'
$ws.Range("B2").Value = 'Unique identification number'
$ws.Range("A3").Value = 'name'
$ws.Range("B3").Value = 'Common name of lake or pool'
$ws.Range("A4").Value = 'area'
$ws.Range("B4").Value = 'Surface area of lake or pool in m2'
$ws.Range("A5").Value = 'elevation'
$ws.Range("B5").Value = 'Elevation of lake or pool in meters above sea level'
$ws.Range("A6").Value = 'latitude'
$ws.Range("B6").Value = 'Latitude of lake or pool centroid in decimal degrees'
$ws.Range("A7").Value = 'longitude'
$ws.Range("B7").Value = 'Longitude of lake or pool centroid in decimal degrees'
$ws.Range("A8").Value = 'label'
$ws.Range("B8").Value = 'Labeling flag'
$ws.Range("C8").Value = '0 or 1'
$ws.Range("A9").Value = 'sample_t2'
$ws.Range("B9").Value = 'Unique water isotope sample name from sampling period 2'
$ws.Range("A10").Value = 'date_t2'
$ws.Range("B10").Value = 'Date water sampled in period 2'
$ws.Range("A11").Value = 'd18O_t2'
$ws.Range("B11").Value = 'The d18O of water sampled in period 2 in per mil'
$ws.Range("A12").Value = 'd2H_t2'
$ws.Range("B12").Value = 'The d2H of water sampled in period 2 in per mil'
$ws.Range("A13").Value = 'dxs_t2'
$ws.Range("B13").Value = 'The dxs of water sampled in period 2 in per mil'
$ws.Range("A14").Value = 'EI_d18O_t2'
$ws.Range("B14").Value = 'The E/I ratio modeled from d18O in period 2'
$ws.Range("A15").Value = 'EI_d2H_t2'
$ws.Range("B15").Value = 'The E/I ratio modeled from d2H in period 2'
$ws.Range("A16").Value = 'EI_dxs_t2'
$ws.Range("B16").Value = 'The E/I ratio modeled from dxs in period 2'
$ws.Range("A17").Value = 'sample_t3'
$ws.Range("B17").Value = 'Unique water isotope sample name from sampling period 3'
$ws.Range("A18").Value = 'date_t3'
$ws.Range("B18").Value = 'Date water sampled in period 3'
$ws.Range("A19").Value = 'd18O_t3'
$ws.Range("B19").Value = 'The d18O of water sampled in period 3 in per mil'
$ws.Range("A20").Value = 'd2H_t3'
$ws.Range("B20").Value = 'The d2H of water sampled in period 3 in per mil'
$ws.Range("A21").Value = 'dxs_t3'
$ws.Range("B21").Value = 'The dxs of water sampled in period 3 in per mil'
$ws.Range("A22").Value = 'EI_d18O_t3'
$ws.Range("B22").Value = 'The E/I ratio modeled from d18O in period 3'
$ws.Range("A23").Value = 'EI_d2H_t3'
$ws.Range("B23").Value = 'The E/I ratio modeled from d2H in period 3'
$ws.Range("A24").Value = 'EI_dxs_t3'
$ws.Range("B24").Value = 'The E/I ratio modeled from dxs in period 3'
$ws.Range("A25").Value = 'laketype'
$ws.Range("B25").Value = 'Category of lake type'
$ws.Range("C25").Value = 'endorheic; headwater; downstream; vale; proglacial; altered'
$ws.Range("A26").Value = 'laketype_num'
$ws.Range("B26").Value = 'Number assigned to category of laketype'
$ws.Range("C26").Value = '1=endorheic; 2=headwater; 3=downstream; 4=vale; 5=proglacial; 6=altered'
$ws.Range("A27").Value = 'surf_area'
$ws.Range("B27").Value = 'Surface area of lake in m2'
$ws.Range("A28").Value = 'basin_name'
$ws.Range("B28").Value = 'Name of drainage basin that lake or pool is located within'
$ws.Range("A29").Value = 'alt_basin'
$ws.Range("B29").Value = 'Alternate name of drainage basin that lake or pool is located within'
$ws.Range("A30").Value = 'lakeshed'
$ws.Range("B30").Value = 'Surface area of lake drainage basin in m2'
$ws.Range("A31").Value = 'dist_gris'
$ws.Range("B31").Value = 'Distance from lake or pool centroid to nearest margin of Greenland ice sheet in m'
$ws.Range("A32").Value = 'dist_ocean'
$ws.Range("B32").Value = 'Distance from lake or pool centroid to nearest ocean coast in m'
$ws.Range("A33").Value = 'main_lakes'
$ws.Range("B33").Value = 'Flag for whether part of the main lakes region'
$ws.Range("C33").Value = '0 or 1'
$ws.Range("A34").Value = 'd18O_infl'
$ws.Range("B34").Value = 'Inferred inflow source water d18O'
$ws.Range("A35").Value = 'd2H_infl'
$ws.Range("B35").Value = 'Inferred inflow source water d2H'
$ws.Range("A36").Value = 'dxs_infl'
$ws.Range("B36").Value = 'Inferred inflow source water dxs'
$ws.Range("A37").Value = 'frze_frac'
$ws.Range("B37").Value = 'Fraction of inflow sourced from frozen season precipitation (Sep-May)'
$ws.Range("A38").Value = 'thaw_frac'
$ws.Range("B38").Value = 'Fraction of inflow sourced from thawed season precipitation (Jun-Aug)'

# Rows 10 and 23 no longer carry a C-column "options" value in the new layout
$ws.Range("C10").ClearContents()
$ws.Range("C23").ClearContents()

# Match the saved view: scrolled down with I30 selected
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("I30").Select()

# --- thule_lakes (sheet11): no data change, just the saved cell selection ---
$ws2 = $wb.Worksheets.Item("thule_lakes")
$ws2.Activate()
$ws2.Range("B4").Select()

# Restore the originally active tab (Overview) so the workbook reopens there
$wb.Worksheets.Item("Overview").Activate()
